$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = '*maa://24633 (54.97), *maa://30515 (69.7), *maa://34787 (72.13), ***maa://20792 (11.93), maa://39402 (85.19), ***maa://29083 (29.41)'
$ws.Range("C4").Value = 'maa://24632 (93.33), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)'
$ws.Range("W11").Value = 'maa://36713 (97.88)'
$ws.Range("W12").Value = 'maa://22753 (91.33), *maa://21485 (76.74), maa://37962 (82.35)'
$ws.Range("C13").Value = 'maa://24999 (91.48), maa://36673 (91.8), maa://25001 (85.51)'
$ws.Range("C14").Value = 'maa://30764 (86.67)'
$ws.Range("AE15").Value = 'maa://21364 (80.68), *maa://22766 (73.0), *maa://36666 (77.27)'
$ws.Range("G17").Value = 'maa://22430 (88.14), maa://39599 (81.82)'
$ws.Range("AA19").Value = '*maa://30709 (61.11), *maa://36668 (52.17)'
$ws.Range("K20").Value = 'maa://41331 (84.44)'
$ws.Range("K23").Value = 'maa://39756 (92.64), maa://39875 (95.92)'
$ws.Range("O23").Value = 'maa://30587 (91.72), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (78.95)'
$ws.Range("W24").Value = 'maa://29988 (86.47), maa://23504 (92.94), **maa://22892 (40.14), *maa://25141 (77.05), maa://36663 (80.7), ***maa://22815 (23.08)'
$ws.Range("AA25").Value = 'maa://31215 (84.34), *maa://24516 (79.07), maa://26001 (87.27)'
$ws.Range("AA26").Value = 'maa://42235 (80.77)'
$ws.Range("C28").Value = 'maa://24465 (90.4), maa://25725 (82.28)'
$ws.Range("S28").Value = 'maa://23263 (94.32), *maa://29765 (60.29)'
$ws.Range("W28").Value = 'maa://39929 (87.14), ***maa://39723 (14.71), maa://41749 (85.0)'
$ws.Range("AE28").Value = 'maa://36660 (93.54), *maa://36701 (64.0)'
$ws.Range("AE29").Value = '*maa://24080 (68.92), ***maa://34960 (8.7)'
$ws.Range("W30").Value = '*maa://39477 (80.0)'
$ws.Range("S32").Value = 'maa://41108 (87.8), maa://41238 (94.74)'
$ws.Range("K35").Value = 'maa://41296 (98.28)'
$ws.Range("O40").Value = 'maa://23278 (95.89), maa://21386 (95.65), maa://36664 (90.7)'
$ws.Range("S44").Value = 'maa://39366 (84.62)'
# F45 changes from the numeric-looking text "4" to "5"; the cell's number
# format is General, so assigning the string directly would be re-typed by
# Excel as a Number (losing the original inlineStr/Text storage). Stage the
# literal text in a scratch cell formatted as Text, then paste-values it in,
# which is how Excel preserves "text that looks like a number" without
# disturbing F45's existing style/number format.
$ws.Range("ZZ1000").NumberFormat = "@"
$ws.Range("ZZ1000").Value = '5'
$ws.Range("ZZ1000").Copy()
$ws.Range("F45").PasteSpecial(-4163)
$ws.Range("ZZ1000").Clear()
$ws.Range("G45").Value = 'maa://21229 (85.47), maa://30807 (95.08), *maa://22767 (52.94), ***maa://20796 (13.79), maa://42459 (100.0)'
$ws.Range("O45").Value = '*maa://36237 (58.33)'
$ws.Range("S45").Value = '**maa://39364 (45.45)'
$ws.Range("G46").Value = 'maa://35931 (92.27)'
$ws.Range("G55").Value = 'maa://32532 (91.96)'
$ws.Range("G58").Value = '*maa://37964 (60.0)'
$excel.CutCopyMode = $false
